# Trade #60 closed at 2026-02-17 08:48:14 - unknown UNKNOWN +0.000%
#
# This script updates the "live trading results" workbook to record the
# closing of trade #60 on the MarketMaking strategy:
#   - Summary sheet: refresh aggregate capital / P&L / trade-count metrics
#   - Strategy Status sheet: refresh the MarketMaking strategy row
#   - All Trades / MarketMaking sheets: append the new closed-trade row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.68   # Current Capital
$summary.Range("B4").Value = -0.32     # Total P&L $
$summary.Range("B5").Value = -0.11     # Total P&L %
$summary.Range("B6").Value = 60        # Total Trades
$summary.Range("B7").Value = 24        # Winning Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking row = row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.68000000000001   # Capital
$status.Range("D4").Value = 60                  # Trades
$status.Range("E4").Value = -0.32               # P&L $
$status.Range("F4").Value = -0.32               # P&L %
$status.Range("G4").Value = 40                  # Win Rate %

# ---------------------------------------------------------------------
# Helper: append the new trade row (#60, row 61) to a trades-log sheet.
# Date/time-looking text must be forced to stay as plain text instead of
# being auto-converted to a date/time serial number by the Excel engine.
# ---------------------------------------------------------------------
function Add-TradeRow60($ws) {
    $row = 61

    $ws.Cells.Item($row, 1).Value = 60        # Trade #

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"   # Date
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "08:48:08"     # Time
    $ws.Cells.Item($row, 3).ClearFormats()

    $ws.Cells.Item($row, 4).Value = "MarketMaking"  # Strategy
    $ws.Cells.Item($row, 5).Value = "UP"            # Side
    $ws.Cells.Item($row, 6).Value = 0.72            # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.78            # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"        # Status
    $ws.Cells.Item($row, 9).Value = 8.333299999999999   # P&L %
    $ws.Cells.Item($row, 10).Value = 0.06               # P&L $
    $ws.Cells.Item($row, 11).Value = 99.68000000000001  # Capital After
    $ws.Cells.Item($row, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"        # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.13                # Duration (min)
}

# ---------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow60 $allTrades

# ---------------------------------------------------------------------
# 4) MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow60 $marketMaking
